$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.467973557154338, 0.99375,              1.479230274545386, 0.9828333333333333, 160.030779, 46.743),
    @(1.467830800771713, 0.9939791666666666,   1.478338003475615, 0.9837083333333334, 159.771541, 46.297),
    @(1.467804515891605, 0.9940625000000001,   1.478516336239821, 0.9834444444444445, 159.500118, 45.947),
    @(1.467977272907893, 0.99390625,           1.478287551948365, 0.9836875,          161.905068, 46.033),
    @(1.46800213470459,  0.9938499999999999,   1.477782581714874, 0.9841333333333333, 161.358039, 46.03),
    @(1.467941167566511, 0.9939305555555555,   1.477706554088187, 0.9842638888888889, 161.117685, 46.693),
    @(1.467813669795082, 0.9940625,            1.477493817954803, 0.9844999999999999, 160.55456,  46.79),
    @(1.467829661448797, 0.9940416666666667,   1.477305639455927, 0.9847083333333333, 161.1254,   47.189),
    @(1.467792299641503, 0.9940879629629631,   1.477287031305033, 0.98475,            159.52721,  47.73),
    @(1.467819216680527, 0.9940645833333335,   1.477159527451434, 0.9848916666666667, 161.770044, 48.745)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
